# Update cryptocurrency price (D) and volume-change (E) cells to the latest
# scraped values, mirroring the GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.847.87'
$ws.Range("E2").Value = '  +0.26%  '

$ws.Range("D3").Value = '2.816.81'
$ws.Range("E3").Value = '  +1.80%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '''354.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.73%  '

$ws.Range("D6").Value = '''111.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.82%  '

$ws.Range("D7").Value = '''0.564'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.46%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("E9").Value = '  +4.56%  '

$ws.Range("D10").Value = '''40.75'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.82%  '

$ws.Range("D11").Value = '''0.0854'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.15%  '

$ws.Range("E12").Value = '  +0.92%  '

$ws.Range("E13").Value = '  +0.02%  '

$ws.Range("D14").Value = '''7.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.99%  '

$ws.Range("D15").Value = '3.261.78'
$ws.Range("E15").Value = '  +1.53%  '

$ws.Range("D16").Value = '2.812.42'
$ws.Range("E16").Value = '  +1.92%  '

$ws.Range("D17").Value = '''0.917'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.72%  '

$ws.Range("D18").Value = '51.767.76'
$ws.Range("E18").Value = '  +0.02%  '

$ws.Range("D19").Value = '''7.56'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.72%  '

$ws.Range("D20").Value = '''3.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.80%  '

$ws.Range("D21").Value = '''13.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.83%  '

$ws.Range("E22").Value = '  +1.54%  '

$ws.Range("D23").Value = '''69.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.03%  '

$ws.Range("D24").Value = '''267.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.82%  '

$ws.Range("D25").Value = '''2.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.70%  '

$ws.Range("D26").Value = '''26.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.93%  '

$ws.Range("E27").Value = '  +0.09%  '

$ws.Range("D28").Value = '''10.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.29%  '

$ws.Range("E29").Value = '  +1.11%  '

$ws.Range("D30").Value = '''0.0479'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +26.00%  '

$ws.Range("D31").Value = '''52.74'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.41%  '

$ws.Range("E32").Value = '  -0.29%  '

$ws.Range("D33").Value = '''34.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.82%  '

$ws.Range("D34").Value = '''5.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.68%  '

$ws.Range("D35").Value = '''5.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.85%  '

$ws.Range("D36").Value = '''0.0844'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.61%  '

$ws.Range("E37").Value = '  -0.17%  '

$ws.Range("D38").Value = '''3.32'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.64%  '

$ws.Range("E39").Value = '  -2.21%  '

$ws.Range("E40").Value = '  -2.61%  '

$ws.Range("E41").Value = '  +0.99%  '

$ws.Range("D42").Value = '''2.54'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.08%  '

$ws.Range("D43").Value = '''23.33'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.41%  '

$ws.Range("D44").Value = '''124.35'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.94%  '

$ws.Range("E45").Value = '  -3.38%  '

$ws.Range("D46").Value = '2.091.98'
$ws.Range("E46").Value = '  +1.36%  '

$ws.Range("E47").Value = '  +1.45%  '

$ws.Range("E48").Value = '  +1.10%  '

$ws.Range("D49").Value = '''5.96'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.94%  '

$ws.Range("D50").Value = '''0.973'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.67%  '

$ws.Range("E51").Value = '  +2.51%  '
